$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the style of the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:50:32.214914",
    "2021-10-05 10:50:32.214926",
    "2021-10-05 10:50:32.214931",
    "2021-10-05 10:50:32.214934",
    "2021-10-05 10:50:32.214937",
    "2021-10-05 10:50:32.214940",
    "2021-10-05 10:50:32.214944",
    "2021-10-05 10:50:32.214947",
    "2021-10-05 10:50:32.214950",
    "2021-10-05 10:50:32.214953",
    "2021-10-05 10:50:32.214956",
    "2021-10-05 10:50:32.214959",
    "2021-10-05 10:50:32.214962",
    "2021-10-05 10:50:32.214965",
    "2021-10-05 10:50:32.214968",
    "2021-10-05 10:50:32.214971",
    "2021-10-05 10:50:32.214974",
    "2021-10-05 10:50:32.214978",
    "2021-10-05 10:50:32.214981",
    "2021-10-05 10:50:32.214984",
    "2021-10-05 10:50:32.214986",
    "2021-10-05 10:50:32.214990",
    "2021-10-05 10:50:32.214992",
    "2021-10-05 10:50:32.214995",
    "2021-10-05 10:50:32.214999",
    "2021-10-05 10:50:32.215002",
    "2021-10-05 10:50:32.215005",
    "2021-10-05 10:50:32.215008",
    "2021-10-05 10:50:32.215011",
    "2021-10-05 10:50:32.215018",
    "2021-10-05 10:50:32.215021",
    "2021-10-05 10:50:32.215024",
    "2021-10-05 10:50:32.215028",
    "2021-10-05 10:50:32.215031",
    "2021-10-05 10:50:32.215034",
    "2021-10-05 10:50:32.215052",
    "2021-10-05 10:50:32.215055",
    "2021-10-05 10:50:32.215058"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
